$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column A / B data (rows 1-60) ----
# Updated dataset for the master-thesis revision: rows 1-40 get new
# computed values and rows 41-60 are newly appended.
$data = New-Object 'object[,]' 60,2
$data[0,0] = 0.027600000000000003
$data[0,1] = 0
$data[1,0] = 0.028601294327848414
$data[1,1] = 0.0018573501969297143
$data[2,0] = 0.029470353337400279
$data[2,1] = 0.0038677611022306132
$data[3,0] = 0.030188157241896781
$data[3,1] = 0.0060305643780009096
$data[4,0] = 0.030733835508605832
$data[4,1] = 0.0083431930172362008
$data[5,0] = 0.031087106735310802
$data[5,1] = 0.01079388752235012
$data[6,0] = 0.031234939598316826
$data[6,1] = 0.013352422527120466
$data[7,0] = 0.031176549100875615
$data[7,1] = 0.015974278671041783
$data[8,0] = 0.030911785928368192
$data[8,1] = 0.018631053788633414
$data[9,0] = 0.03044297038281529
$data[9,1] = 0.021298681609308586
$data[10,0] = 0.029773757926402866
$data[10,1] = 0.023957023192411407
$data[11,0] = 0.028906578979824633
$data[11,1] = 0.026591095841910946
$data[12,0] = 0.02784029336220735
$data[12,1] = 0.029190915449779486
$data[13,0] = 0.026568387157913348
$data[13,1] = 0.031750288247939198
$data[14,0] = 0.02508232959574519
$data[14,1] = 0.034261625626509841
$data[15,0] = 0.023387971244823155
$data[15,1] = 0.036705054500755022
$data[16,0] = 0.021496069034282091
$data[16,1] = 0.039060298834845873
$data[17,0] = 0.01941781558147522
$data[17,1] = 0.041310044771103359
$data[18,0] = 0.017163718790668255
$data[18,1] = 0.043439789110732582
$data[19,0] = 0.014743255892444252
$data[19,1] = 0.045437163357475691
$data[20,0] = 0.012165412443227119
$data[20,1] = 0.047291085457566037
$data[21,0] = 0.0094401317148825915
$data[21,1] = 0.048991083679273455
$data[22,0] = 0.0065731265071155582
$data[22,1] = 0.050528095608198792
$data[23,0] = 0.0035587101735440065
$data[23,1] = 0.051893504589564057
$data[24,0] = 0.00038487777691255073
$data[24,1] = 0.05307552762253475
$data[25,0] = -0.0029656632378995736
$data[25,1] = 0.054057172135719887
$data[26,0] = -0.0065137780112569134
$data[26,1] = 0.054814329294629406
$data[27,0] = -0.010280391163594629
$data[27,1] = 0.055314322445300344
$data[28,0] = -0.014242098992306954
$data[28,1] = 0.055525649606486593
$data[29,0] = -0.018349660980057997
$data[29,1] = 0.055426106263442368
$data[30,0] = -0.022541691954509561
$data[30,1] = 0.055006520803706466
$data[31,0] = -0.026750911017145065
$data[31,1] = 0.054273101884406379
$data[32,0] = -0.030901619919662168
$data[32,1] = 0.053252793953515079
$data[33,0] = -0.034954033205832125
$data[33,1] = 0.051969498911221539
$data[34,0] = -0.038920636153665783
$data[34,1] = 0.050417200840244626
$data[35,0] = -0.042837807450057899
$data[35,1] = 0.048559065524243168
$data[36,0] = -0.046754051893864466
$data[36,1] = 0.046321954660580068
$data[37,0] = -0.050739121783802153
$data[37,1] = 0.043566780473718515
$data[38,0] = -0.054877922507997522
$data[38,1] = 0.040050570250862075
$data[39,0] = -0.059142549923437843
$data[39,1] = 0.035540945240013846
$data[40,0] = -0.029330544062522997
$data[40,1] = 0.047314086321386312
$data[41,0] = -0.032529482723521859
$data[41,1] = 0.046036423933391177
$data[42,0] = -0.035700487754425136
$data[42,1] = 0.044526346086814342
$data[43,0] = -0.038826879500871027
$data[43,1] = 0.042780594296847503
$data[44,0] = -0.041893077029166143
$data[44,1] = 0.040795173331206143
$data[45,0] = -0.044888361294267327
$data[45,1] = 0.038559902643114818
$data[46,0] = -0.047801814062131499
$data[46,1] = 0.036060468616549247
$data[47,0] = -0.050627641032028874
$data[47,1] = 0.033267661361455744
$data[48,0] = -0.053346791264648759
$data[48,1] = 0.030158323008008834
$data[49,0] = -0.055914392174097043
$data[49,1] = 0.026748445565299591
$data[50,0] = -0.058267480416534186
$data[50,1] = 0.023114456023117409
$data[51,0] = -0.060374419001357552
$data[51,1] = 0.019307262676770937
$data[52,0] = -0.062228189217733688
$data[52,1] = 0.01532780449043274
$data[53,0] = -0.063816171932330246
$data[53,1] = 0.01119009996470678
$data[54,0] = -0.065124935683038704
$data[54,1] = 0.0069197293738260694
$data[55,0] = -0.066144721707631096
$data[55,1] = 0.0025334184544201989
$data[56,0] = -0.066864416563743229
$data[56,1] = -0.0020197733689228318
$data[57,0] = -0.067251466756522327
$data[57,1] = -0.0068221642828780911
$data[58,0] = -0.067249074379770735
$data[58,1] = -0.011925622520987033
$data[59,0] = -0.066807520966436773
$data[59,1] = -0.017255582931883569

$ws.Range("A1:B60").Value2 = $data

# ---- Column A width ----
$ws.Columns.Item(1).ColumnWidth = 14.8
